$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.542.37"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.493.47"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.66"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.12"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.79"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.872.49"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.48"
$ws.Range("E15").Value = "  +6.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.476.74"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.761"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.599.18"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.91"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.18"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.93"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.61"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.46"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.14"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.26"
$ws.Range("E33").Value = "  +6.30%  "
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.48"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.99"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.52"
$ws.Range("E43").Value = "  -8.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.952.51"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.79"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.727.84"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.23"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.21"
$ws.Range("E51").Value = "  -2.72%  "
